$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing text formatting for the Price (D) / Volume (E) columns while
# writing the refreshed values; some of the new figures look like plain numbers
# (e.g. "232.83"), and without forcing text format Excel would coerce them to
# numeric cells instead of the inline/shared strings the sheet already uses.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '37.771.07'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '2.076.73'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '232.83'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '57.21'
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = '0.387'
$ws.Range("E9").Value = '  +1.17%  '
$ws.Range("D10").Value = '0.0787'
$ws.Range("E10").Value = '  +3.08%  '
$ws.Range("E11").Value = '  +2.79%  '
$ws.Range("D12").Value = '2.369.95'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '14.44'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").Value = '20.97'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = '0.761'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '5.27'
$ws.Range("E16").Value = '  +2.23%  '
$ws.Range("D17").Value = '2.073.83'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '37.663.76'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("E19").Value = '  -3.76%  '
$ws.Range("D20").Value = '70.59'
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("D21").Value = '0.0₃0821'
$ws.Range("E21").Value = '  +1.00%  '
$ws.Range("D22").Value = '228.01'
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("D26").Value = '170.22'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("E27").Value = '  +10.08%  '
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("D30").Value = '19.38'
$ws.Range("E30").Value = '  +1.94%  '
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("D32").Value = '4.63'
$ws.Range("E32").Value = '  +3.68%  '
$ws.Range("D33").Value = '0.0625'
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("D34").Value = '4.62'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = '2.51'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("E36").Value = '  +3.46%  '
$ws.Range("D37").Value = '3.40'
$ws.Range("E37").Value = '  +5.28%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  -4.18%  '
$ws.Range("D40").Value = '0.100'
$ws.Range("E40").Value = '  +7.32%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").Value = '97.97'
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("D44").Value = '1.447.84'
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("D46").Value = '4.13'
$ws.Range("E46").Value = '  -4.44%  '
$ws.Range("E47").Value = '  +2.81%  '
$ws.Range("D48").Value = '15.60'
$ws.Range("E48").Value = '  +2.75%  '
$ws.Range("E49").Value = '  +3.66%  '
$ws.Range("D50").Value = '3.00'
$ws.Range("E50").Value = '  +1.23%  '
$ws.Range("D51").Value = '2.266.55'
$ws.Range("E51").Value = '  +0.39%  '
# Restore the default (unstyled) cell style so no visible formatting changes
# are introduced - only the underlying values change.
$dataRange.Style = "Normal"
